$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$ws = $wb.Worksheets.Item("Rushing")
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 17
$ws.Range("F2").Value = 8

$ws.Range("C4").Value = 218
$ws.Range("D4").Value = 104
$ws.Range("F4").Value = 44

$ws.Range("C9").Value = 7

# --- Receiving sheet ---
$ws = $wb.Worksheets.Item("Receiving")
$ws.Range("C2").Value = 58
$ws.Range("D2").Value = 50

$ws.Range("C3").Value = 35
$ws.Range("D3").Value = 29
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 6

$ws.Range("C5").Value = 90
$ws.Range("D5").Value = 67
$ws.Range("G5").Value = 10

$ws.Range("C6").Value = 107
$ws.Range("D6").Value = 77
$ws.Range("E6").Value = 37
$ws.Range("F6").Value = 22
$ws.Range("G6").Value = 12

$ws.Range("C7").Value = 105
$ws.Range("D7").Value = 75
$ws.Range("E7").Value = 47
$ws.Range("F7").Value = 25
$ws.Range("G7").Value = 16
$ws.Range("H7").Value = 10

$ws.Range("C11").Value = 69

$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 11
